$d = $word.ActiveDocument

$replacements = @(
    @("624×6=3744", "833×3=2499"),
    @("352×5=1760", "470×2=940"),
    @("314×6=1884", "941×8=7528"),
    @("622×7=4354", "978×2=1956"),
    @("410×8=3280", "179×8=1432"),
    @("837×9=7533", "101×6=606"),
    @("707×7=4949", "828×5=4140"),
    @("150×5=750", "514×9=4626"),
    @("852×4=3408", "647×7=4529"),
    @("871×8=6968", "924×3=2772"),
    @("616×3=1848", "870×7=6090"),
    @("142×4=568", "626×5=3130"),
    @("379×2=758", "179×4=716"),
    @("858×7=6006", "677×5=3385"),
    @("395×4=1580", "949×8=7592"),
    @("773×7=5411", "975×8=7800"),
    @("828×9=7452", "355×7=2485"),
    @("137×4=548", "770×2=1540"),
    @("187×3=561", "249×3=747"),
    @("259×6=1554", "567×9=5103"),
    @("700×3=2100", "822×3=2466"),
    @("647×5=3235", "421×9=3789"),
    @("300×2=600", "627×3=1881"),
    @("672×4=2688", "797×8=6376"),
    @("152×2=304", "836×5=4180")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
